$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Commit message: "Updated excel documents to read ANep as well as removed
# any data sets that did not include pT > 1."
# In this sheet, the "AN" label (column J, rows 2-13) should read "ANep".
$ws.Range("J2:J13").Value = "ANep"

# Move the active selection to G14, matching the saved cursor position
# recorded in the updated sheetView.
$ws.Range("G14").Select()
